$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq "Administrator, Miss Dina Nasr") {
        $cell.Value2 = "Miss Dina Nasr, Administrator"
    }
}
